$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" text in cell A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 14.81 = 62666.67 pesos`n✅ 62666.67 pesos = 14.81 = 980.06 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- tasas: update the rate values in N10, O10, N12, O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 67.5
$wsTasas.Range("O10").Value = 4230
$wsTasas.Range("N12").Value = 4231
$wsTasas.Range("O12").Value = 66.17
